# Update countries & provincias Spain
# - New country "Congo" enters the ranked list (sorted desc by "Casos totales"),
#   shifting the countries below it down by one row; "Puerto Rico" and "Zambia"
#   also swap order. Refreshed case/death/recovered counters for several countries
#   below, and the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados..." timestamp: 11:22 -> 11:52
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 5 de Abril de 2020 a las 11:52"

# Refreshed statistics for existing rows (Espana, Alemania, Suiza, Suecia, Rumania,
# Panama, Sri Lanka, Brunei, Isla de Man)
$ws.Cells.Item(5, 2).Value = 130759
$ws.Cells.Item(5, 3).Value = 4591
$ws.Cells.Item(5, 4).Value = 38080
$ws.Cells.Item(5, 5).Value = 80261
$ws.Cells.Item(5, 7).Value = 471
$ws.Cells.Item(5, 8).Value = 12418

$ws.Cells.Item(7, 2).Value = 96108
$ws.Cells.Item(7, 3).Value = 16
$ws.Cells.Item(7, 5).Value = 68264

$ws.Cells.Item(13, 2).Value = 20510
$ws.Cells.Item(13, 3).Value = 5
$ws.Cells.Item(13, 5).Value = 13429

$ws.Cells.Item(22, 6).Value = 520

$ws.Cells.Item(32, 5).Value = 3136
$ws.Cells.Item(32, 7).Value = 2
$ws.Cells.Item(32, 8).Value = 148

$ws.Cells.Item(44, 6).Value = 75

$ws.Cells.Item(110, 2).Value = 167
$ws.Cells.Item(110, 3).Value = 1
$ws.Cells.Item(110, 5).Value = 133

$ws.Cells.Item(117, 4).Value = 73
$ws.Cells.Item(117, 5).Value = 61

$ws.Cells.Item(120, 2).Value = 127
$ws.Cells.Item(120, 3).Value = 1
$ws.Cells.Item(120, 5).Value = 126

# "Congo" newly enters the table at row 140; rows 140-155 shift down one place
# (with "Puerto Rico"/"Zambia" swapping order) and pick up refreshed counts
$ws.Cells.Item(140, 1).Value = "Congo"
$ws.Cells.Item(140, 2).Value = 45
$ws.Cells.Item(140, 3).Value = 23
$ws.Cells.Item(140, 4).Value = 2
$ws.Cells.Item(140, 5).Value = 38
$ws.Cells.Item(140, 7).Value = 3
$ws.Cells.Item(140, 8).Value = 5

$ws.Cells.Item(141, 1).Value = "Macao"
$ws.Cells.Item(141, 2).Value = 44
$ws.Cells.Item(141, 4).Value = 10
$ws.Cells.Item(141, 5).Value = 34
$ws.Cells.Item(141, 8).Value = 0

$ws.Cells.Item(142, 1).Value = "Mali"
$ws.Cells.Item(142, 4).Value = 1
$ws.Cells.Item(142, 5).Value = 37

$ws.Cells.Item(143, 1).Value = "Togo"
$ws.Cells.Item(143, 2).Value = 41
$ws.Cells.Item(143, 4).Value = 17
$ws.Cells.Item(143, 5).Value = 21
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 8).Value = 3

$ws.Cells.Item(144, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(144, 2).Value = 40
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 5).Value = 40
$ws.Cells.Item(144, 6).Value = 1
$ws.Cells.Item(144, 8).Value = 0

$ws.Cells.Item(146, 1).Value = "Zambia"
$ws.Cells.Item(146, 2).Value = 39
$ws.Cells.Item(146, 4).Value = 2
$ws.Cells.Item(146, 5).Value = 36
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 8).Value = 1

$ws.Cells.Item(147, 1).Value = "Etiopia"
$ws.Cells.Item(147, 2).Value = 38
$ws.Cells.Item(147, 4).Value = 4
$ws.Cells.Item(147, 5).Value = 34
$ws.Cells.Item(147, 6).Value = 1

$ws.Cells.Item(148, 1).Value = "Bermudas"
$ws.Cells.Item(148, 2).Value = 37
$ws.Cells.Item(148, 4).Value = 14
$ws.Cells.Item(148, 5).Value = 23
$ws.Cells.Item(148, 8).Value = 0

$ws.Cells.Item(149, 1).Value = "Islas Caimanes"
$ws.Cells.Item(149, 2).Value = 35
$ws.Cells.Item(149, 4).Value = 1
$ws.Cells.Item(149, 5).Value = 33

$ws.Cells.Item(150, 1).Value = "Guam"
$ws.Cells.Item(150, 2).Value = 32
$ws.Cells.Item(150, 5).Value = 31
$ws.Cells.Item(150, 8).Value = 1

$ws.Cells.Item(151, 1).Value = "Eritrea"
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 5).Value = 29
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 8).Value = 0

$ws.Cells.Item(152, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(152, 2).Value = 29
$ws.Cells.Item(152, 3).Value = 5
$ws.Cells.Item(152, 4).Value = 7
$ws.Cells.Item(152, 5).Value = 20
$ws.Cells.Item(152, 6).Value = 6
$ws.Cells.Item(152, 8).Value = 2

$ws.Cells.Item(153, 1).Value = "Bahamas"
$ws.Cells.Item(153, 2).Value = 28
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 5).Value = 24
$ws.Cells.Item(153, 6).Value = 1

$ws.Cells.Item(154, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(154, 2).Value = 25
$ws.Cells.Item(154, 3).Value = 2
$ws.Cells.Item(154, 4).Value = 6
$ws.Cells.Item(154, 5).Value = 15

$ws.Cells.Item(155, 1).Value = "Guyana"
$ws.Cells.Item(155, 2).Value = 24
$ws.Cells.Item(155, 3).Value = 1
$ws.Cells.Item(155, 4).Value = 0
$ws.Cells.Item(155, 5).Value = 20
$ws.Cells.Item(155, 8).Value = 4
